$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values. Each price cell in column D is stored
# as text (not a number) in the source sheet, e.g. "243.59", so we force
# Text formatting before the write to stop Excel from auto-coercing the
# numeric-looking string into a real number (which would also silently
# drop meaningful trailing zeros, e.g. "0.03010" -> 0.0301). ClearFormats()
# afterwards restores the cell to its original (unstyled) appearance so only
# the value itself changes.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "243.59"
Set-TextValue $ws.Range("D3") "23.24"
Set-TextValue $ws.Range("D4") "5.644"
Set-TextValue $ws.Range("D5") "0.05851"
Set-TextValue $ws.Range("D7") "6.477"
Set-TextValue $ws.Range("D8") "1.318"
Set-TextValue $ws.Range("D9") "0.7984"
Set-TextValue $ws.Range("D10") "0.1459"
Set-TextValue $ws.Range("D11") "0.07627"
Set-TextValue $ws.Range("D12") "0.03254"
Set-TextValue $ws.Range("D13") "0.03010"
Set-TextValue $ws.Range("D14") "0.09241"
Set-TextValue $ws.Range("D15") "0.001659"
Set-TextValue $ws.Range("D16") "3.416"
Set-TextValue $ws.Range("D17") "0.04755"
Set-TextValue $ws.Range("D18") "0.0005996"
Set-TextValue $ws.Range("D19") "0.006222"
Set-TextValue $ws.Range("D20") "0.001069"
Set-TextValue $ws.Range("D21") "0.003830"
Set-TextValue $ws.Range("D24") "2.209"
Set-TextValue $ws.Range("D25") "0.3337"
Set-TextValue $ws.Range("D26") "0.1252"
Set-TextValue $ws.Range("D27") "0.0004004"
Set-TextValue $ws.Range("E27") "26UpBotsUBXTWorstin24h"
Set-TextValue $ws.Range("D40") "0.04312"
Set-TextValue $ws.Range("D41") "0.007099"
Set-TextValue $ws.Range("D42") "0.1055"
Set-TextValue $ws.Range("D43") "0.003244"
Set-TextValue $ws.Range("D44") "0.008725"
Set-TextValue $ws.Range("E44") "43LocalTradersLCT"
Set-TextValue $ws.Range("D46") "0.00005754"
Set-TextValue $ws.Range("D48") "0.7861"
Set-TextValue $ws.Range("D49") "0.1013"
